# binary classification model rebuild with hyperparameter tuning
# Update the evaluation metrics on the active sheet to reflect the newly
# retrained / hyperparameter-tuned model's reported scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - class "0"
$ws.Range("B4").Value = 0.9826994189117803
$ws.Range("C4").Value = 0.8130463286713286
$ws.Range("D4").Value = 0.8898588854341066

# Row 5 - class "1"
$ws.Range("B5").Value = 0.0300453514739229
$ws.Range("C5").Value = 0.2880434782608696
$ws.Range("D5").Value = 0.05441478439425051

# Row 6 - accuracy
$ws.Range("B6").Value = 0.8026992287917738
$ws.Range("C6").Value = 0.8026992287917738
$ws.Range("D6").Value = 0.8026992287917738
$ws.Range("E6").Value = 0.8026992287917738

# Row 7 - macro avg
$ws.Range("B7").Value = 0.5063723851928515
$ws.Range("C7").Value = 0.5505449034660991
$ws.Range("D7").Value = 0.4721368349141786

# Row 8 - weighted avg
$ws.Range("B8").Value = 0.9639238888765868
$ws.Range("C8").Value = 0.8026992287917738
$ws.Range("D8").Value = 0.8733934061505447

# Row 9 - overall
$ws.Range("F9").Value = 0.5659587949414715
$ws.Range("G9").Value = 0.2943228781223297
